# progression analysis complex query
#
# Adds a small "Variante" (V1..V4) lookup table to the interactive-short
# sheet (columns K:N, mirroring the existing Query/Post/Comment/Profile
# table in B3:F3), tidies up the left-over placeholder cells around it,
# and leaves the workbook with the "interactive-short" tab active.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("interactive-short")

# The old placeholder cells in J1:N5 carried no data, only stray
# (mostly default) formatting left over from an earlier layout. Clear the
# ones that aren't part of the new Variante block so they stop being
# serialized as empty cells - but leave K1 / J3 / K3 alone, they keep
# their existing (pre-existing) formatting in the new layout too.
$ws.Range("J1").Clear()
$ws.Range("L1").Clear()
$ws.Range("M1").Clear()
$ws.Range("N1").Clear()
$ws.Range("J2").Clear()
$ws.Range("L3").Clear()
$ws.Range("M3").Clear()
$ws.Range("N3").Clear()
$ws.Range("J4").Clear()
$ws.Range("L4").Clear()
$ws.Range("M4").Clear()
$ws.Range("N4").Clear()
$ws.Range("J5").Clear()
$ws.Range("L5").Clear()
$ws.Range("M5").Clear()
$ws.Range("N5").Clear()

# New "Variante" mini-table: header row + four variant labels (V1..V4),
# same shape as the Query table header (Query/Post/Comment/Profile).
$ws.Range("K2").Value = "Variante"
$ws.Range("L2").Value = "Post"
$ws.Range("M2").Value = "Comment"
$ws.Range("N2").Value = "Profile"

$ws.Range("K3").Value = "V1"
$ws.Range("K4").Value = "V2"
$ws.Range("K5").Value = "V3"
$ws.Range("K6").Value = "V4"

# The K column used to be sized for older, wider content; re-fit it now
# that it only holds "Variante"/"V1".."V4".
$ws.Range("K1").ColumnWidth = 8.3

# interactive-discover column F widened slightly (cosmetic).
$wsDiscover = $wb.Worksheets.Item("interactive-discover")
$wsDiscover.Range("F1").ColumnWidth = 24.3

# Make "interactive-short" the active tab (was "interactive-discover").
$ws.Activate()
